$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the rows/columns that held the old (now removed) data so no stale
# values remain beyond the new, smaller table (old table was A1:B6).
$ws.Range("A1:B6").ClearContents()

# Header row
$ws.Range("A1").Value = "Variable"
$ws.Range("B1").Value = "Importance"

# Data rows (relabeled variables / updated importance values)
$ws.Range("A2").Value = "HH Private Car Ownership (%)"
$ws.Range("B2").Value = 12.19851365203616

$ws.Range("A3").Value = "Population Density (/sq. km)"
$ws.Range("B3").Value = 9.231980484501491

$ws.Range("A4").Value = "GDP per capita (2018 US$)"
$ws.Range("B4").Value = 3.033333333333333
